$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to rounded (2 decimal place) values
$ws.Range("B5").Value = 6.73
$ws.Range("C5").Value = 4.72
$ws.Range("D5").Value = 0.78
$ws.Range("E5").Value = 14.32
$ws.Range("F5").Value = 11.85
$ws.Range("G5").Value = 5.29
$ws.Range("H5").Value = 25.51
$ws.Range("I5").Value = 8.140000000000001
$ws.Range("J5").Value = 3.54
$ws.Range("K5").Value = 5.26
$ws.Range("L5").Value = 5.84
$ws.Range("M5").Value = 6.01
$ws.Range("N5").Value = 1.7
$ws.Range("O5").Value = 5.26
$ws.Range("P5").Value = 7.45
$ws.Range("Q5").Value = 4.6
$ws.Range("R5").Value = 0.75
$ws.Range("S5").Value = 0.42
$ws.Range("T5").Value = 72.83
$ws.Range("U5").Value = 14.96
$ws.Range("V5").Value = 4.86
$ws.Range("W5").Value = 9.890000000000001
$ws.Range("X5").Value = 5.33
$ws.Range("Y5").Value = 0.58
$ws.Range("Z5").Value = 11.77
$ws.Range("AA5").Value = 4.29
$ws.Range("AB5").Value = 3.92
$ws.Range("AC5").Value = 4.58
$ws.Range("AD5").Value = 6.06
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 23.37
$ws.Range("AG5").Value = 2.66
$ws.Range("AH5").Value = 6.08

# Delete row 6 entirely (shifts dimension back to A1:AH5)
$ws.Range("A6:AH6").Delete()
